$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("AC1").Value = "wnb-调节6Hz_20161230_113123_ASIC_EEG"
$ws.Range("AD1").Value = "wnb-调节6Hz_20170110_113300_ASIC_EEG"
$ws.Range("AC2").Value = 0.864951768488746
$ws.Range("AD2").Value = 0.95469255663430419
$ws.Range("AC3").Value = 0.86956521739130432
$ws.Range("AD3").Value = 0.9078498293515358
$ws.Range("A1:AD3").Select() | Out-Null
